# InstructionChecklist.xlsx — "Added tests and implemented instructions"
#
# The checklist tracks, per ARM instruction (rows 2-29), whether it has an
# entry in the "InstructionsInterfaces" column (E). This commit marks every
# instruction from MOV through B as done ("y") while leaving the
# not-yet-covered instructions (BL, DCD, EQU, FILL, END, rows 30-34) and the
# trailing blank rows (35-42) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark column E ("InstructionsInterfaces") as done for rows 2 through 29.
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 5).Value = "y"
}

# Restore the view: scroll so row 10 is at the top-left, and leave the
# active selection on E24 (matches the saved workbook view state).
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("E24").Select()
